# "Generate Report for Handback"
#
# For each language sheet (zh-cn, de-de) this stamps rows 2 & 3 with the
# handback report info:
#   - "Latest Target File"   (col F) -> same source .md file as col A
#   - "Latest Handback File" (col G) -> the .xlf file that was handed back
#     (same file referenced by col D, the "Latest Handoff File")
#   - "Latest Handback DateTime" (col H) -> the real timestamp the handback
#     completed at, replacing the "0001-01-01 00:00:00" placeholder.
# It also flips every "Status" cell (Overview + both language sheets) from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdFile   = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.md"
$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/c5c0ad920ee49906314e121359f41c1f8eae48a7/e2e/cd3d2eed-9657-46d0-a4ae-8a7a25503f11.md"

$zhFile   = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.zh-cn.xlf"
$zhUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e3339dc276efe42aaaee1be2d92d0bc8f5d81ef/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.zh-cn.xlf"

$deFile   = "cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.de-de.xlf"
$deUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a4521c3229c03c060ae107c6306187839f5b167/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cd3d2eed-9657-46d0-a4ae-8a7a25503f11.72c9c9f2dede7c896e651960276eab8bc9184b97.de-de.xlf"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: flip every Status cell (B/C, rows 2-3) -----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) on rows 2-3
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Target File (F) / Latest Handback File (G) hyperlinks, rows 2-3
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhUrl, "", "", $zhFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhUrl, "", "", $zhFile)

# Latest Handback DateTime (H) now has a real timestamp instead of the
# "0001-01-01 00:00:00" placeholder.
$wsZh.Range("H2").Value = "2016-03-23 09:15:08"
$wsZh.Range("H3").Value = "2016-03-23 09:15:08"

# --- de-de sheet --------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (C) on rows 2-3
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Target File (F) / Latest Handback File (G) hyperlinks, rows 2-3
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deUrl, "", "", $deFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deUrl, "", "", $deFile)

# Latest Handback DateTime (H) -- de-de's handback completed a few seconds
# after zh-cn's, so it gets its own (later) timestamp.
$wsDe.Range("H2").Value = "2016-03-23 09:15:14"
$wsDe.Range("H3").Value = "2016-03-23 09:15:14"
